$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply cell updates from the cryptos list refresh (values scraped on Mon Jul 31 03:35:26 UTC 2023)
# Each target cell is forced to Text format before assignment so that numeric-looking
# strings (e.g. '1.001', '29.450.79') are preserved as text exactly as in the source data,
# then the style is reset to Normal so no stray formatting is introduced.
$updates = @(
    @{ Cell = 'D2'; Value = '29.450.79' }
    @{ Cell = 'E2'; Value = '  +0.31%  ' }
    @{ Cell = 'D3'; Value = '1.870.86' }
    @{ Cell = 'E4'; Value = '  -0.13%  ' }
    @{ Cell = 'D5'; Value = '243.92' }
    @{ Cell = 'E5'; Value = '  +0.44%  ' }
    @{ Cell = 'D6'; Value = '0.7063' }
    @{ Cell = 'E6'; Value = '  -0.71%  ' }
    @{ Cell = 'E7'; Value = '  -0.10%  ' }
    @{ Cell = 'D8'; Value = '0.3155' }
    @{ Cell = 'E8'; Value = '  +0.57%  ' }
    @{ Cell = 'D9'; Value = '0.07882' }
    @{ Cell = 'E9'; Value = '  -1.69%  ' }
    @{ Cell = 'D10'; Value = '24.68' }
    @{ Cell = 'E10'; Value = '  -1.77%  ' }
    @{ Cell = 'D11'; Value = '0.08014' }
    @{ Cell = 'E11'; Value = '  -3.79%  ' }
    @{ Cell = 'D12'; Value = '1.894.69' }
    @{ Cell = 'E12'; Value = '  +0.66%  ' }
    @{ Cell = 'D14'; Value = '94.25' }
    @{ Cell = 'E14'; Value = '  -0.69%  ' }
    @{ Cell = 'D15'; Value = '0.7058' }
    @{ Cell = 'E15'; Value = '  -1.83%  ' }
    @{ Cell = 'D16'; Value = '6.494' }
    @{ Cell = 'E16'; Value = '  +2.11%  ' }
    @{ Cell = 'B17'; Value = 'WrappedBTC' }
    @{ Cell = 'C17'; Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc' }
    @{ Cell = 'D17'; Value = '29.502.79' }
    @{ Cell = 'E17'; Value = '  +0.30%  ' }
    @{ Cell = 'B18'; Value = 'ShibaInu' }
    @{ Cell = 'C18'; Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib' }
    @{ Cell = 'D18'; Value = '0.000008374' }
    @{ Cell = 'E18'; Value = '  -3.67%  ' }
    @{ Cell = 'D19'; Value = '256.93' }
    @{ Cell = 'E19'; Value = '  +5.57%  ' }
    @{ Cell = 'D20'; Value = '2.139.23' }
    @{ Cell = 'E20'; Value = '  -1.77%  ' }
    @{ Cell = 'E21'; Value = '  -1.00%  ' }
    @{ Cell = 'E22'; Value = '  -0.09%  ' }
    @{ Cell = 'D23'; Value = '7.638' }
    @{ Cell = 'E23'; Value = '  -2.85%  ' }
    @{ Cell = 'D24'; Value = '1.001' }
    @{ Cell = 'E24'; Value = '  -0.13%  ' }
    @{ Cell = 'D25'; Value = '0.1557' }
    @{ Cell = 'E25'; Value = '  -1.10%  ' }
    @{ Cell = 'D26'; Value = '9.076' }
    @{ Cell = 'E26'; Value = '  -0.19%  ' }
    @{ Cell = 'D27'; Value = '161.03' }
    @{ Cell = 'E27'; Value = '  -1.56%  ' }
    @{ Cell = 'E28'; Value = '  +1.16%  ' }
    @{ Cell = 'D29'; Value = '1.503' }
    @{ Cell = 'E29'; Value = '  -0.32%  ' }
    @{ Cell = 'D30'; Value = '4.342' }
    @{ Cell = 'D31'; Value = '4.256' }
    @{ Cell = 'E31'; Value = '  -2.39%  ' }
    @{ Cell = 'D32'; Value = '1.209' }
    @{ Cell = 'E32'; Value = '  +0.21%  ' }
    @{ Cell = 'D33'; Value = '0.05325' }
    @{ Cell = 'D34'; Value = '1.899' }
    @{ Cell = 'E34'; Value = '  -2.27%  ' }
    @{ Cell = 'B35'; Value = 'ARBITRUM' }
    @{ Cell = 'C35'; Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb' }
    @{ Cell = 'D35'; Value = '1.174' }
    @{ Cell = 'E35'; Value = '  -0.56%  ' }
    @{ Cell = 'B36'; Value = 'ImmutableX' }
    @{ Cell = 'C36'; Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx' }
    @{ Cell = 'D36'; Value = '0.7481' }
    @{ Cell = 'E36'; Value = '  -3.87%  ' }
    @{ Cell = 'D37'; Value = '2.715' }
    @{ Cell = 'E37'; Value = '  +0.88%  ' }
    @{ Cell = 'D38'; Value = '0.01881' }
    @{ Cell = 'E38'; Value = '  -0.41%  ' }
    @{ Cell = 'D39'; Value = '1.263.60' }
    @{ Cell = 'E39'; Value = '  -0.54%  ' }
    @{ Cell = 'E40'; Value = '  +0.21%  ' }
    @{ Cell = 'D41'; Value = '0.9003' }
    @{ Cell = 'E41'; Value = '  -2.14%  ' }
    @{ Cell = 'D42'; Value = '108.86' }
    @{ Cell = 'E42'; Value = '  -4.39%  ' }
    @{ Cell = 'D43'; Value = '71.84' }
    @{ Cell = 'D44'; Value = '5.952' }
    @{ Cell = 'E44'; Value = '  -8.99%  ' }
    @{ Cell = 'B45'; Value = 'PaxDollar' }
    @{ Cell = 'C45'; Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp' }
    @{ Cell = 'D45'; Value = '1.000' }
    @{ Cell = 'E45'; Value = '  -0.13%  ' }
    @{ Cell = 'B46'; Value = 'BabyDogeCoin' }
    @{ Cell = 'C46'; Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge' }
    @{ Cell = 'D46'; Value = '0.00000000130' }
    @{ Cell = 'E46'; Value = '  +1.75%  ' }
    @{ Cell = 'D47'; Value = '2.038.59' }
    @{ Cell = 'E47'; Value = '  +0.35%  ' }
    @{ Cell = 'D48'; Value = '1.811' }
    @{ Cell = 'E48'; Value = '  -0.02%  ' }
    @{ Cell = 'E49'; Value = '  -0.53%  ' }
    @{ Cell = 'E50'; Value = '  -0.50%  ' }
    @{ Cell = 'D51'; Value = '0.4326' }
    @{ Cell = 'E51'; Value = '  -1.34%  ' }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    $rng.NumberFormat = "@"
    $rng.Value = $u.Value
    $rng.Style = "Normal"
}
